$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates derived from the diff (prices / 1h volume % / two row swaps).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "54.172.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.270.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.04%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "497.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  -0.87%  "
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.336"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.69"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.670.73"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.69"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "54.140.50"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.65%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.271.68"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.72%  "
$ws.Range("E18").Value = "  +1.84%  "
$ws.Range("E19").Value = "  +2.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "302.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.45%  "
$ws.Range("E21").Value = "  -2.61%  "
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.44%  "
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.149"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "170.47"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("E28").Value = "  +0.27%  "
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.91"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.68%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0684"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.45%  "
$ws.Range("E31").Value = "  +0.33%  "
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.72"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.75%  "
$ws.Range("B34").Value = "SuiNetwork"
$ws.Range("C34").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.956"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +10.66%  "
$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("E36").Value = "  -1.19%  "
$ws.Range("E37").Value = "  +1.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.372"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.84%  "
$ws.Range("E39").Value = "  -0.16%  "
$ws.Range("E40").Value = "  +0.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.81"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "124.82"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.02%  "
$ws.Range("E43").Value = "  +2.10%  "
$ws.Range("E44").Value = "  -0.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.543"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.77%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "238.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.39%  "
$ws.Range("E47").Value = "  -0.67%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0204"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.75"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.37%  "
$ws.Range("E50").Value = "  -1.06%  "
$ws.Range("E51").Value = "  -0.29%  "
